$wb = $excel.ActiveWorkbook

$wsUsers = $wb.Worksheets.Item("UserList")
$wsRequests = $wb.Worksheets.Item("RequestList")

# Populate the RequestList sheet with header row + one sample data row
$headers = @("username", "First Name", "Last Name", "DoB", "Card #", "Last Accessed", "Employee Status", "Password")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsRequests.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$wsRequests.Cells.Item(2, 1).Value = "banana"
$wsRequests.Cells.Item(2, 2).Value = "Chris"
$wsRequests.Cells.Item(2, 3).Value = "Moticska"

# Update selections on each sheet
$wsUsers.Range("A1:H1").Select()
$wsRequests.Range("F3").Select()

# Make RequestList the active sheet/tab
$wsRequests.Activate()
